$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15755.97779219088
$ws.Range("C2").Value = 863358.8077491384
$ws.Range("B3").Value = 23775.99974660412
$ws.Range("C3").Value = 1302821.003241439
$ws.Range("B4").Value = 38780.97304745853
$ws.Range("C4").Value = 2125028.043019956
